$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column Q - copy formatting from P1 (same header style), then set value
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q1").Value = "n_studies"
$excel.CutCopyMode = $false

# Fill Q2:Q17 with the number of studies (36) for each replication row
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 17).Value = 36
}
